$d = $word.ActiveDocument

# The document currently ends with a single empty paragraph
# (ListParagraph style, ilvl 0, numId 1). We turn that paragraph into
# the new "Polish List" heading, insert four new list items after it,
# and finish with a new "Bug List" heading that becomes the document's
# new final paragraph.

$lastPara = $d.Paragraphs.Last
$target = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:contextualSpacing/>
    <w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr><w:t>Polish List</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Player</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Camera</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Middle mouse reset camera position behind player</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Movement</w:t></w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@
[void]$target.InsertXML($xml)

# The document again ends with a lone empty paragraph (same style as
# before). Turn it into the "Bug List" heading. To get Word to accept
# the new paragraph-mark formatting (pStyle/numPr/contextualSpacing/
# rPr) on what is the very last paragraph of the body, we first give
# ourselves a fresh trailing paragraph to merge into: type the text,
# add one more paragraph after it, then InsertXML across both so the
# new formatting lands on the merged (now final) paragraph mark.

$bugListPara = $d.Paragraphs.Last
$bugRange = $bugListPara.Range
$bugRange.Text = "Bug List"
$bugRange.InsertParagraphAfter()

$bugListPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$trailingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bugTarget = $d.Range($bugListPara.Range.Start, $trailingPara.Range.End)

$bugXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:contextualSpacing/>
    <w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr><w:t>Bug List</w:t></w:r>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@
[void]$bugTarget.InsertXML($bugXml)
